$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26 in the original sheet held a stray "Международные стандарты" value in
# column B, with column A blank, sitting between row 25 ("7. Сопоставимость с
# международными данными / стандартами") and row 27 ("8. Ссылки и
# документация"). Fold that orphaned B26 value up into B25 (which belongs
# with the "Сопоставимость..." label), then remove the now-empty row so
# everything below shifts up by one.
$ws.Range("B25").Value = $ws.Range("B26").Text
$ws.Rows("26:26").Delete()

# Reflect the post-edit selection state (cursor parked just past the new
# last data row, with the view scrolled back to the top of the sheet).
[void]$ws.Range("A33").Select()
